# "Generate Report for handoff"
#
# The localization-status report is regenerated: the row for the
# "5828addd-6cf5-4598-b31f-b644ff7b492e" source file moves from the 2nd
# data row to the last data row (alphabetical-ish reorder against the
# other two source files), and its status flips from
# "Handed back: in sync with en-US" to "Ready for handoff" because it has
# just been queued for a fresh handoff (its handoff datetime is bumped
# forward as well). The other two rows simply shift up one position.
#
# This touches all three worksheets (Overview, zh-cn, de-de) plus the
# hyperlinks that decorate their file-name / xlf-name columns.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "ffff3a438823-8e45-4206-ad5a-b92aae305d78.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "Handed back: in sync with en-US"

$ws.Range("A3").Value = "ffffff40bbcbd8-a44f-4250-bc3c-15d7c526b25d.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "Handed back: in sync with en-US"

$ws.Range("A4").Value = "5828addd-6cf5-4598-b31f-b644ff7b492e.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "Ready for handoff"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7781df4ec36f55449ca1d69da3383f13d6929925/e2e/5828addd-6cf5-4598-b31f-b644ff7b492e.md", [Type]::Missing, [Type]::Missing, "ffff3a438823-8e45-4206-ad5a-b92aae305d78.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/7781df4ec36f55449ca1d69da3383f13d6929925/e2e/ffff3a438823-8e45-4206-ad5a-b92aae305d78.md", [Type]::Missing, [Type]::Missing, "ffffff40bbcbd8-a44f-4250-bc3c-15d7c526b25d.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/7781df4ec36f55449ca1d69da3383f13d6929925/e2e/ffffff40bbcbd8-a44f-4250-bc3c-15d7c526b25d.md", [Type]::Missing, [Type]::Missing, "5828addd-6cf5-4598-b31f-b644ff7b492e.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/7781df4ec36f55449ca1d69da3383f13d6929925/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# 2. zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "ffff3a438823-8e45-4206-ad5a-b92aae305d78.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "ca3c0fb3-b850-4ae1-9914-d36efa6e197e.c6a2a43d5404e307aba0548c33fcaa8e348f522d.zh-cn.xlf"
$ws.Range("D2").Value = "2016-01-22 02:57:45"
$ws.Range("E2").Value = "ca3c0fb3-b850-4ae1-9914-d36efa6e197e.md"
$ws.Range("F2").Value = "ca3c0fb3-b850-4ae1-9914-d36efa6e197e.c6a2a43d5404e307aba0548c33fcaa8e348f522d.zh-cn.xlf"
$ws.Range("G2").Value = "2016-01-22 02:58:33"
$ws.Range("H2").Value = "Include"

$ws.Range("A3").Value = "ffffff40bbcbd8-a44f-4250-bc3c-15d7c526b25d.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "ca3c0fb3-b850-4ae1-9914-d36efa6e197e.c6a2a43d5404e307aba0548c33fcaa8e348f522d.zh-cn.xlf"
$ws.Range("D3").Value = "2016-01-22 02:57:45"
$ws.Range("E3").Value = "ca3c0fb3-b850-4ae1-9914-d36efa6e197e.md"
$ws.Range("F3").Value = "ca3c0fb3-b850-4ae1-9914-d36efa6e197e.c6a2a43d5404e307aba0548c33fcaa8e348f522d.zh-cn.xlf"
$ws.Range("G3").Value = "2016-01-22 02:58:33"
$ws.Range("H3").Value = "Include"

$ws.Range("A4").Value = "5828addd-6cf5-4598-b31f-b644ff7b492e.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "5828addd-6cf5-4598-b31f-b644ff7b492e.072c592b6be096debc01a3b10d2c576139b6327b.zh-cn.xlf"
$ws.Range("D4").Value = "2016-01-22 03:04:38"
$ws.Range("E4").Value = "5828addd-6cf5-4598-b31f-b644ff7b492e.md"
$ws.Range("F4").Value = "5828addd-6cf5-4598-b31f-b644ff7b492e.072c592b6be096debc01a3b10d2c576139b6327b.zh-cn.xlf"
$ws.Range("G4").Value = "2016-01-22 03:03:25"
$ws.Range("H4").Value = "Include"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7781df4ec36f55449ca1d69da3383f13d6929925/e2e/5828addd-6cf5-4598-b31f-b644ff7b492e.md", [Type]::Missing, [Type]::Missing, "ffff3a438823-8e45-4206-ad5a-b92aae305d78.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d00b2cbf221e954b491d8bb318ad6b765c163135/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/5828addd-6cf5-4598-b31f-b644ff7b492e.072c592b6be096debc01a3b10d2c576139b6327b.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "ca3c0fb3-b850-4ae1-9914-d36efa6e197e.c6a2a43d5404e307aba0548c33fcaa8e348f522d.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/90b6bdd01669bd403fd1551c8872cd30b55f7abb/e2e/5828addd-6cf5-4598-b31f-b644ff7b492e.md", [Type]::Missing, [Type]::Missing, "ca3c0fb3-b850-4ae1-9914-d36efa6e197e.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/bde63bed76480870edcbb31895dbdf3b7fe0a329/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/5828addd-6cf5-4598-b31f-b644ff7b492e.072c592b6be096debc01a3b10d2c576139b6327b.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "ca3c0fb3-b850-4ae1-9914-d36efa6e197e.c6a2a43d5404e307aba0548c33fcaa8e348f522d.zh-cn.xlf") | Out-Null

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/7781df4ec36f55449ca1d69da3383f13d6929925/e2e/ffff3a438823-8e45-4206-ad5a-b92aae305d78.md", [Type]::Missing, [Type]::Missing, "ffffff40bbcbd8-a44f-4250-bc3c-15d7c526b25d.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ef777b5f9c49a32dfbea3c35718b6f0c01350e61/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ca3c0fb3-b850-4ae1-9914-d36efa6e197e.c6a2a43d5404e307aba0548c33fcaa8e348f522d.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "ca3c0fb3-b850-4ae1-9914-d36efa6e197e.c6a2a43d5404e307aba0548c33fcaa8e348f522d.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/3b7b4e262c6cec5639734583ee7aeb94812f1623/e2e/ca3c0fb3-b850-4ae1-9914-d36efa6e197e.md", [Type]::Missing, [Type]::Missing, "ca3c0fb3-b850-4ae1-9914-d36efa6e197e.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/46263d875c6225855df8eca6af19a21a5a12bb94/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ca3c0fb3-b850-4ae1-9914-d36efa6e197e.c6a2a43d5404e307aba0548c33fcaa8e348f522d.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "ca3c0fb3-b850-4ae1-9914-d36efa6e197e.c6a2a43d5404e307aba0548c33fcaa8e348f522d.zh-cn.xlf") | Out-Null

$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/7781df4ec36f55449ca1d69da3383f13d6929925/e2e/ffffff40bbcbd8-a44f-4250-bc3c-15d7c526b25d.md", [Type]::Missing, [Type]::Missing, "5828addd-6cf5-4598-b31f-b644ff7b492e.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ef777b5f9c49a32dfbea3c35718b6f0c01350e61/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ca3c0fb3-b850-4ae1-9914-d36efa6e197e.c6a2a43d5404e307aba0548c33fcaa8e348f522d.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "5828addd-6cf5-4598-b31f-b644ff7b492e.072c592b6be096debc01a3b10d2c576139b6327b.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/3b7b4e262c6cec5639734583ee7aeb94812f1623/e2e/ca3c0fb3-b850-4ae1-9914-d36efa6e197e.md", [Type]::Missing, [Type]::Missing, "5828addd-6cf5-4598-b31f-b644ff7b492e.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/46263d875c6225855df8eca6af19a21a5a12bb94/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ca3c0fb3-b850-4ae1-9914-d36efa6e197e.c6a2a43d5404e307aba0548c33fcaa8e348f522d.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "5828addd-6cf5-4598-b31f-b644ff7b492e.072c592b6be096debc01a3b10d2c576139b6327b.zh-cn.xlf") | Out-Null

$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/7781df4ec36f55449ca1d69da3383f13d6929925/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# 3. de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "ffff3a438823-8e45-4206-ad5a-b92aae305d78.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "ca3c0fb3-b850-4ae1-9914-d36efa6e197e.c6a2a43d5404e307aba0548c33fcaa8e348f522d.de-de.xlf"
$ws.Range("D2").Value = "2016-01-22 02:57:59"
$ws.Range("E2").Value = "ca3c0fb3-b850-4ae1-9914-d36efa6e197e.md"
$ws.Range("F2").Value = "ca3c0fb3-b850-4ae1-9914-d36efa6e197e.c6a2a43d5404e307aba0548c33fcaa8e348f522d.de-de.xlf"
$ws.Range("G2").Value = "2016-01-22 02:58:57"
$ws.Range("H2").Value = "Include"

$ws.Range("A3").Value = "ffffff40bbcbd8-a44f-4250-bc3c-15d7c526b25d.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "ca3c0fb3-b850-4ae1-9914-d36efa6e197e.c6a2a43d5404e307aba0548c33fcaa8e348f522d.de-de.xlf"
$ws.Range("D3").Value = "2016-01-22 02:57:59"
$ws.Range("E3").Value = "ca3c0fb3-b850-4ae1-9914-d36efa6e197e.md"
$ws.Range("F3").Value = "ca3c0fb3-b850-4ae1-9914-d36efa6e197e.c6a2a43d5404e307aba0548c33fcaa8e348f522d.de-de.xlf"
$ws.Range("G3").Value = "2016-01-22 02:58:57"
$ws.Range("H3").Value = "Include"

$ws.Range("A4").Value = "5828addd-6cf5-4598-b31f-b644ff7b492e.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "5828addd-6cf5-4598-b31f-b644ff7b492e.072c592b6be096debc01a3b10d2c576139b6327b.de-de.xlf"
$ws.Range("D4").Value = "2016-01-22 03:04:53"
$ws.Range("E4").Value = "5828addd-6cf5-4598-b31f-b644ff7b492e.md"
$ws.Range("F4").Value = "5828addd-6cf5-4598-b31f-b644ff7b492e.072c592b6be096debc01a3b10d2c576139b6327b.de-de.xlf"
$ws.Range("G4").Value = "2016-01-22 03:03:49"
$ws.Range("H4").Value = "Include"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7781df4ec36f55449ca1d69da3383f13d6929925/e2e/5828addd-6cf5-4598-b31f-b644ff7b492e.md", [Type]::Missing, [Type]::Missing, "ffff3a438823-8e45-4206-ad5a-b92aae305d78.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c7f54c9ca12c34391f67aeb1ca43113702abca2d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/5828addd-6cf5-4598-b31f-b644ff7b492e.072c592b6be096debc01a3b10d2c576139b6327b.de-de.xlf", [Type]::Missing, [Type]::Missing, "ca3c0fb3-b850-4ae1-9914-d36efa6e197e.c6a2a43d5404e307aba0548c33fcaa8e348f522d.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/09e15498730f89a57a07db1827edeed7c9f8aefa/e2e/5828addd-6cf5-4598-b31f-b644ff7b492e.md", [Type]::Missing, [Type]::Missing, "ca3c0fb3-b850-4ae1-9914-d36efa6e197e.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c8ef6a2432d1e345b012e80e7643d2b8146a2a48/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/5828addd-6cf5-4598-b31f-b644ff7b492e.072c592b6be096debc01a3b10d2c576139b6327b.de-de.xlf", [Type]::Missing, [Type]::Missing, "ca3c0fb3-b850-4ae1-9914-d36efa6e197e.c6a2a43d5404e307aba0548c33fcaa8e348f522d.de-de.xlf") | Out-Null

$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/7781df4ec36f55449ca1d69da3383f13d6929925/e2e/ffff3a438823-8e45-4206-ad5a-b92aae305d78.md", [Type]::Missing, [Type]::Missing, "ffffff40bbcbd8-a44f-4250-bc3c-15d7c526b25d.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c71040fbf23ced8ac95e4d6232df08703a279021/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ca3c0fb3-b850-4ae1-9914-d36efa6e197e.c6a2a43d5404e307aba0548c33fcaa8e348f522d.de-de.xlf", [Type]::Missing, [Type]::Missing, "ca3c0fb3-b850-4ae1-9914-d36efa6e197e.c6a2a43d5404e307aba0548c33fcaa8e348f522d.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/76f2836c197a47e56b6ebd7d081484de0f2523d6/e2e/ca3c0fb3-b850-4ae1-9914-d36efa6e197e.md", [Type]::Missing, [Type]::Missing, "ca3c0fb3-b850-4ae1-9914-d36efa6e197e.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b87b2d8d43f225ab58014b445d025f90ecf4e26e/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ca3c0fb3-b850-4ae1-9914-d36efa6e197e.c6a2a43d5404e307aba0548c33fcaa8e348f522d.de-de.xlf", [Type]::Missing, [Type]::Missing, "ca3c0fb3-b850-4ae1-9914-d36efa6e197e.c6a2a43d5404e307aba0548c33fcaa8e348f522d.de-de.xlf") | Out-Null

$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/7781df4ec36f55449ca1d69da3383f13d6929925/e2e/ffffff40bbcbd8-a44f-4250-bc3c-15d7c526b25d.md", [Type]::Missing, [Type]::Missing, "5828addd-6cf5-4598-b31f-b644ff7b492e.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c71040fbf23ced8ac95e4d6232df08703a279021/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ca3c0fb3-b850-4ae1-9914-d36efa6e197e.c6a2a43d5404e307aba0548c33fcaa8e348f522d.de-de.xlf", [Type]::Missing, [Type]::Missing, "5828addd-6cf5-4598-b31f-b644ff7b492e.072c592b6be096debc01a3b10d2c576139b6327b.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/76f2836c197a47e56b6ebd7d081484de0f2523d6/e2e/ca3c0fb3-b850-4ae1-9914-d36efa6e197e.md", [Type]::Missing, [Type]::Missing, "5828addd-6cf5-4598-b31f-b644ff7b492e.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b87b2d8d43f225ab58014b445d025f90ecf4e26e/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ca3c0fb3-b850-4ae1-9914-d36efa6e197e.c6a2a43d5404e307aba0548c33fcaa8e348f522d.de-de.xlf", [Type]::Missing, [Type]::Missing, "5828addd-6cf5-4598-b31f-b644ff7b492e.072c592b6be096debc01a3b10d2c576139b6327b.de-de.xlf") | Out-Null

$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/7781df4ec36f55449ca1d69da3383f13d6929925/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

$wb.Worksheets.Item("Overview").Activate()
